$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new "Save" column, formatted like the other header cells
$ws.Range("H1").Value = "Save"
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = -4108
$ws.Range("H1").VerticalAlignment = -4160
$ws.Range("H1").Borders.LineStyle = 1

# Save flag values for rows 2-38: 1 when sum (col G) >= 10, else 0
$saveValues = @(0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1,0,1,0,0,0,1,0,0,1,0,0,1,0,0,1,0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
